$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 331
$ws.Range("F2").Value = 45108
$ws.Range("G2").Value = 29952
$ws.Range("H2").Value = 45139
$ws.Range("E3").Value = 29891
$ws.Range("F3").Value = 45078
$ws.Range("E4").Value = 29891
$ws.Range("F4").Value = 45078
$ws.Range("G4").Value = 29952
$ws.Range("H4").Value = 45139
$ws.Range("E5").Value = 29921
$ws.Range("F5").Value = 45108
$ws.Range("G5").Value = 29952
$ws.Range("H5").Value = 45139
$ws.Range("C6").Value = 439
$ws.Range("F6").Value = 45078
$ws.Range("G6").Value = 29952
$ws.Range("H6").Value = 45139
$ws.Range("E7").Value = 29891
$ws.Range("F7").Value = 45078
$ws.Range("G7").Value = 29952
$ws.Range("H7").Value = 45139
$ws.Range("D8").Value = 406
$ws.Range("E8").Value = 29891
$ws.Range("F8").Value = 45078
$ws.Range("H8").Value = 45139
$ws.Range("E9").Value = 29891
$ws.Range("F9").Value = 45078
$ws.Range("G9").Value = 29952
$ws.Range("H9").Value = 45139
$ws.Range("C10").Value = 498
$ws.Range("D10").Value = 479
$ws.Range("F10").Value = 45078
$ws.Range("H10").Value = 45139
$ws.Range("E11").Value = 29860
$ws.Range("F11").Value = 45047
$ws.Range("G11").Value = 29952
$ws.Range("H11").Value = 45139
$ws.Range("C12").Value = 367
$ws.Range("D12").Value = 348
$ws.Range("F12").Value = 45078
$ws.Range("H12").Value = 45139
$ws.Range("C13").Value = 463
$ws.Range("F13").Value = 45078
$ws.Range("G13").Value = 29952
$ws.Range("H13").Value = 45139
$ws.Range("C14").Value = 420
$ws.Range("D14").Value = 392
$ws.Range("F14").Value = 45078
$ws.Range("H14").Value = 45139
$ws.Range("C15").Value = 381
$ws.Range("F15").Value = 45108
$ws.Range("G15").Value = 29921
$ws.Range("H15").Value = 45139
$ws.Range("C16").Value = 451
$ws.Range("D16").Value = 406
$ws.Range("F16").Value = 45078
$ws.Range("H16").Value = 45139
$ws.Range("C17").Value = 366
$ws.Range("D17").Value = 390
$ws.Range("F17").Value = 45078
$ws.Range("H17").Value = 45139
$ws.Range("E18").Value = 29891
$ws.Range("F18").Value = 45078
$ws.Range("G18").Value = 29952
$ws.Range("H18").Value = 45139
$ws.Range("D19").Value = 394
$ws.Range("E19").Value = 29891
$ws.Range("F19").Value = 45078
$ws.Range("H19").Value = 45139
$ws.Range("C20").Value = 477
$ws.Range("F20").Value = 45078
$ws.Range("G20").Value = 29952
$ws.Range("H20").Value = 45139
$ws.Range("E21").Value = 29891
$ws.Range("F21").Value = 45078
$ws.Range("G21").Value = 29952
$ws.Range("H21").Value = 45139
$ws.Range("D22").Value = 376
$ws.Range("E22").Value = 29860
$ws.Range("F22").Value = 45047
$ws.Range("H22").Value = 45139
$ws.Range("E23").Value = 29891
$ws.Range("F23").Value = 45078
$ws.Range("G23").Value = 29952
$ws.Range("H23").Value = 45139
$ws.Range("C24").Value = 390
$ws.Range("F24").Value = 45078
$ws.Range("G24").Value = 29952
$ws.Range("H24").Value = 45139
$ws.Range("E25").Value = 29891
$ws.Range("F25").Value = 45078
$ws.Range("G25").Value = 29952
$ws.Range("H25").Value = 45139
$ws.Range("C26").Value = 451
$ws.Range("D26").Value = 394
$ws.Range("F26").Value = 45108
$ws.Range("H26").Value = 45139
$ws.Range("C27").Value = 499
$ws.Range("D27").Value = 406
$ws.Range("F27").Value = 45108
$ws.Range("H27").Value = 45139
$ws.Range("C28").Value = 355
$ws.Range("D28").Value = 363
$ws.Range("F28").Value = 45078
$ws.Range("H28").Value = 45139
